$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -80.8442
$ws.Range("B2").Value = -80.7652

$ws.Range("A3").Value = 32.1791
$ws.Range("B3").Value = 32.2458

$ws.Range("A4").Value = -80.051
$ws.Range("B4").Value = -80.1309

$ws.Range("A5").Value = 32.8437
$ws.Range("B5").Value = 32.7773
